$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 309.76923
$ws.Range("I9").Value = 308.55554
$ws.Range("J9").Value = 312.5
$ws.Range("K9").Value = 308.55554
$ws.Range("L9").Value = 312.5
$ws.Range("M9").Value = -139.55554
$ws.Range("N9").Value = -650.5
$ws.Range("H18").Value = 1889.5
$ws.Range("I18").Value = 1629
$ws.Range("K18").Value = 1629
$ws.Range("M18").Value = -1345
$ws.Range("H28").Value = 949.4706
$ws.Range("I28").Value = 617.4286
$ws.Range("K28").Value = 617.4286
$ws.Range("M28").Value = -132.4286
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 1000
$ws.Range("K32").Value = 1000
$ws.Range("M32").Value = -674
$ws.Range("H64").Value = 3133
$ws.Range("I64").Value = 2999
$ws.Range("J64").Value = 3200
$ws.Range("K64").Value = 2999
$ws.Range("L64").Value = 3200
$ws.Range("M64").Value = -2751
$ws.Range("N64").Value = -3696
$ws.Range("H67").Value = 3133
$ws.Range("I67").Value = 2999
$ws.Range("J67").Value = 3200
$ws.Range("K67").Value = 2999
$ws.Range("L67").Value = 3200
$ws.Range("M67").Value = -2141
$ws.Range("N67").Value = -4916
$ws.Range("H111").Value = 3032
$ws.Range("I111").Value = 1299.25
$ws.Range("J111").Value = 6497.5
$ws.Range("K111").Value = 3897.75
$ws.Range("L111").Value = 19492.5
$ws.Range("M111").Value = -830.75
$ws.Range("N111").Value = -25626.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 45999.5
$ws.Range("J43").Value = 45999
$ws.Range("L43").Value = 45999
$ws.Range("N43").Value = -46625
$ws.Range("H61").Value = 2449.6
$ws.Range("I61").Value = 2449.6
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2449.6
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2237.6
$ws.Range("N61").ClearContents()
$ws.Range("H97").Value = 1193.2
$ws.Range("I97").Value = 1214.7778
$ws.Range("J97").Value = 999
$ws.Range("K97").Value = 1214.7778
$ws.Range("L97").Value = 999
$ws.Range("M97").Value = -718.7778000000001
$ws.Range("N97").Value = -1991
$ws.Range("H123").Value = 80000
$ws.Range("J123").Value = 80000
$ws.Range("L123").Value = 80000
$ws.Range("N123").Value = -89800
$ws.Range("H132").Value = 1811.5
$ws.Range("I132").Value = 1772.2667
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 5316.800099999999
$ws.Range("L132").Value = 7200
$ws.Range("M132").Value = -2786.800099999999
$ws.Range("N132").Value = -12260
$ws.Range("H136").Value = 2449.6
$ws.Range("I136").Value = 2449.6
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7348.799999999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4798.799999999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 15349.75
$ws.Range("I26").Value = 15349.75
$ws.Range("K26").Value = 15349.75
$ws.Range("M26").Value = -15057.75
$ws.Range("H82").Value = 51666.332
$ws.Range("I82").Value = 27500
$ws.Range("J82").Value = 99999
$ws.Range("K82").Value = 27500
$ws.Range("L82").Value = 99999
$ws.Range("M82").Value = -27117
$ws.Range("N82").Value = -100765
$ws.Range("H85").Value = 51666.332
$ws.Range("I85").Value = 27500
$ws.Range("J85").Value = 99999
$ws.Range("K85").Value = 27500
$ws.Range("L85").Value = 99999
$ws.Range("M85").Value = -26174
$ws.Range("N85").Value = -102651
$ws.Range("H96").Value = 9279
$ws.Range("I96").Value = 9279
$ws.Range("K96").Value = 9279
$ws.Range("M96").Value = -6533
$ws.Range("H134").Value = 2286.5908
$ws.Range("I134").Value = 2114.2222
$ws.Range("J134").Value = 3062.25
$ws.Range("K134").Value = 6342.6666
$ws.Range("L134").Value = 9186.75
$ws.Range("M134").Value = -3807.6666
$ws.Range("N134").Value = -14256.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 500
$ws.Range("I17").Value = 500
$ws.Range("K17").Value = 500
$ws.Range("M17").Value = -326
$ws.Range("H31").Value = 4275.647
$ws.Range("I31").Value = 3242.5
$ws.Range("J31").Value = 6755.2
$ws.Range("K31").Value = 3242.5
$ws.Range("L31").Value = 6755.2
$ws.Range("M31").Value = -2947.5
$ws.Range("N31").Value = -7345.2
$ws.Range("H34").Value = 4275.647
$ws.Range("I34").Value = 3242.5
$ws.Range("J34").Value = 6755.2
$ws.Range("K34").Value = 3242.5
$ws.Range("L34").Value = 6755.2
$ws.Range("M34").Value = -3040.5
$ws.Range("N34").Value = -7159.2
$ws.Range("H52").Value = 89000
$ws.Range("J52").Value = 89000
$ws.Range("L52").Value = 89000
$ws.Range("N52").Value = -89588
$ws.Range("H58").Value = 2723.4092
$ws.Range("I58").Value = 1575.1666
$ws.Range("J58").Value = 4101.3
$ws.Range("K58").Value = 1575.1666
$ws.Range("L58").Value = 4101.3
$ws.Range("M58").Value = -1372.1666
$ws.Range("N58").Value = -4507.3
$ws.Range("H132").Value = 5999.5
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470
$ws.Range("H134").Value = 2060.9312
$ws.Range("I134").Value = 1742
$ws.Range("J134").Value = 3591.8
$ws.Range("K134").Value = 5226
$ws.Range("L134").Value = 10775.4
$ws.Range("M134").Value = -2691
$ws.Range("N134").Value = -15845.4
$ws.Range("H136").Value = 2723.4092
$ws.Range("I136").Value = 1575.1666
$ws.Range("J136").Value = 4101.3
$ws.Range("K136").Value = 4725.4998
$ws.Range("L136").Value = 12303.9
$ws.Range("M136").Value = -2175.4998
$ws.Range("N136").Value = -17403.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 161
$ws.Range("I103").Value = 178.75
$ws.Range("J103").Value = 90
$ws.Range("K103").Value = 536.25
$ws.Range("L103").Value = 270
$ws.Range("M103").Value = 342.75
$ws.Range("N103").Value = -2028
$ws.Range("H128").Value = 1127177
$ws.Range("I128").Value = 1127177
$ws.Range("K128").Value = 3381531
$ws.Range("M128").Value = -3376551
$ws.Range("H131").Value = 1113.1052
$ws.Range("I131").Value = 525.6667
$ws.Range("K131").Value = 1577.0001
$ws.Range("M131").Value = 3462.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8719.6
$ws.Range("I80").Value = 3299
$ws.Range("J80").Value = 12333.333
$ws.Range("K80").Value = 3299
$ws.Range("L80").Value = 12333.333
$ws.Range("M80").Value = -2301
$ws.Range("N80").Value = -14329.333
$ws.Range("H83").Value = 8719.6
$ws.Range("I83").Value = 3299
$ws.Range("J83").Value = 12333.333
$ws.Range("K83").Value = 16495
$ws.Range("L83").Value = 61666.665
$ws.Range("M83").Value = -11503
$ws.Range("N83").Value = -71650.66500000001
$ws.Range("H107").Value = 90
$ws.Range("I107").Value = 90
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 90
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1830
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 2127.4324
$ws.Range("I132").Value = 1894.0646
$ws.Range("J132").Value = 3333.1667
$ws.Range("K132").Value = 5682.1938
$ws.Range("L132").Value = 9999.500100000001
$ws.Range("M132").Value = -3152.1938
$ws.Range("N132").Value = -15059.5001
$ws.Range("H137").Value = 60704
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 60704
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 60704
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -70904

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 12319.786
$ws.Range("I16").Value = 11329.667
$ws.Range("J16").Value = 13062.375
$ws.Range("K16").Value = 11329.667
$ws.Range("L16").Value = 13062.375
$ws.Range("M16").Value = -11159.667
$ws.Range("N16").Value = -13402.375
$ws.Range("H22").Value = 6190.8237
$ws.Range("I22").Value = 4840.25
$ws.Range("J22").Value = 7391.3335
$ws.Range("K22").Value = 4840.25
$ws.Range("L22").Value = 7391.3335
$ws.Range("M22").Value = -4545.25
$ws.Range("N22").Value = -7981.3335
$ws.Range("H27").Value = 6190.8237
$ws.Range("I27").Value = 4840.25
$ws.Range("J27").Value = 7391.3335
$ws.Range("K27").Value = 4840.25
$ws.Range("L27").Value = 7391.3335
$ws.Range("M27").Value = -4733.25
$ws.Range("N27").Value = -7605.3335
$ws.Range("H34").Value = 4997.5
$ws.Range("I34").Value = 4997.5
$ws.Range("K34").Value = 4997.5
$ws.Range("M34").Value = -4825.5
$ws.Range("H46").Value = 2931.6333
$ws.Range("I46").Value = 1937.4667
$ws.Range("J46").Value = 3925.8
$ws.Range("K46").Value = 1937.4667
$ws.Range("L46").Value = 3925.8
$ws.Range("M46").Value = -1749.4667
$ws.Range("N46").Value = -4301.8
$ws.Range("H93").Value = 2450
$ws.Range("I93").Value = 2450
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2450
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -1202
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 5225.8667
$ws.Range("I132").Value = 4486.75
$ws.Range("K132").Value = 13460.25
$ws.Range("M132").Value = -10930.25
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 69999.5
$ws.Range("J80").Value = 69999.5
$ws.Range("L80").Value = 69999.5
$ws.Range("N80").Value = -71995.5
$ws.Range("H81").Value = 10702.277
$ws.Range("I81").Value = 7410.25
$ws.Range("J81").Value = 11642.857
$ws.Range("K81").Value = 14820.5
$ws.Range("L81").Value = 23285.714
$ws.Range("M81").Value = -13759.5
$ws.Range("N81").Value = -25407.714
$ws.Range("H83").Value = 69999.5
$ws.Range("J83").Value = 69999.5
$ws.Range("L83").Value = 209998.5
$ws.Range("N83").Value = -219982.5
$ws.Range("H84").Value = 10702.277
$ws.Range("I84").Value = 7410.25
$ws.Range("J84").Value = 11642.857
$ws.Range("K84").Value = 74102.5
$ws.Range("L84").Value = 116428.57
$ws.Range("M84").Value = -68798.5
$ws.Range("N84").Value = -127036.57
$ws.Range("H96").Value = 2250
$ws.Range("I96").Value = 1500
$ws.Range("K96").Value = 1500
$ws.Range("M96").Value = -127
$ws.Range("H132").Value = 95373.7
$ws.Range("I132").Value = 118592.125
$ws.Range("K132").Value = 355776.375
$ws.Range("M132").Value = -353246.375
$ws.Range("H136").Value = 1606.4445
$ws.Range("I136").Value = 1475.9231
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 4427.7693
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -1877.7693
$ws.Range("N136").Value = -20100
